# Generate Report for Handoff
#
# Refreshes the localization-status report for a new handoff pass:
#   - zh-cn / de-de status moves from "Handed back: in sync with en-US" to
#     "Ready for handoff" on the Overview sheet and on each language sheet.
#   - The associated "Latest Handoff Datetime" / "Latest HO Xliff Generate
#     Date" timestamps are bumped to the new handoff run.
#   - The "Status" column on each sheet is re-sized to fit the new, shorter
#     status text (it was sized for the old, longer message).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value = "Ready for handoff"       # Status column
$dede.Range("C2").Value = "Ready for handoff"       # Status column

# --- Refresh handoff timestamps for this run ---
$zhcn.Range("H2").Value = "2016-08-21 01:04:20"     # zh-cn Latest Handoff Datetime
$dede.Range("H2").Value = "2016-08-21 01:04:25"     # de-de Latest Handoff Datetime
$overview.Range("G2").Value = "2016-08-21 01:04:25" # Latest HO Xliff Generate Date

# --- Re-fit the "Status" columns to the new, shorter status text ---
$overview.Range("E:F").ColumnWidth = 16.33
$zhcn.Range("C:C").ColumnWidth = 16.33
$dede.Range("C:C").ColumnWidth = 16.33
